# Apply cell value updates per the scheduled runner diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 749.2143
$ws.Range("J55").Value = 1625
$ws.Range("L55").Value = 1625
$ws.Range("N55").Value = -2053
$ws.Range("H76").Value = 13099.956
$ws.Range("J76").Value = 13033
$ws.Range("L76").Value = 13033
$ws.Range("N76").Value = -13663
$ws.Range("H79").Value = 13099.956
$ws.Range("J79").Value = 13033
$ws.Range("L79").Value = 13033
$ws.Range("N79").Value = -15217
$ws.Range("H82").Value = 3221.2222
$ws.Range("J82").Value = 20999
$ws.Range("L82").Value = 62997
$ws.Range("N82").Value = -63809
$ws.Range("H85").Value = 3221.2222
$ws.Range("J85").Value = 20999
$ws.Range("L85").Value = 62997
$ws.Range("N85").Value = -65805
$ws.Range("H88").Value = 5920.4443
$ws.Range("I88").Value = 1947.5
$ws.Range("J88").Value = 7055.5713
$ws.Range("K88").Value = 1947.5
$ws.Range("L88").Value = 7055.5713
$ws.Range("M88").Value = -1541.5
$ws.Range("N88").Value = -7867.5713
$ws.Range("H91").Value = 5920.4443
$ws.Range("I91").Value = 1947.5
$ws.Range("J91").Value = 7055.5713
$ws.Range("K91").Value = 1947.5
$ws.Range("L91").Value = 7055.5713
$ws.Range("M91").Value = -543.5
$ws.Range("N91").Value = -9863.5713
$ws.Range("H98").Value = 1842.6923
$ws.Range("I98").Value = 1544
$ws.Range("J98").Value = 2406.889
$ws.Range("K98").Value = 1544
$ws.Range("L98").Value = 2406.889
$ws.Range("M98").Value = -46
$ws.Range("N98").Value = -5402.889
$ws.Range("H107").Value = 1408.125
$ws.Range("J107").Value = 2612.25
$ws.Range("L107").Value = 2612.25
$ws.Range("N107").Value = -6452.25
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("H122").Value = 1842.6923
$ws.Range("I122").Value = 1544
$ws.Range("J122").Value = 2406.889
$ws.Range("K122").Value = 4632
$ws.Range("L122").Value = 7220.667
$ws.Range("M122").Value = -2182
$ws.Range("N122").Value = -12120.667
$ws.Range("H132").Value = 3601.3044
$ws.Range("I132").Value = 3310.9092
$ws.Range("J132").Value = 9990
$ws.Range("K132").Value = 9932.7276
$ws.Range("L132").Value = 29970
$ws.Range("M132").Value = -7402.7276
$ws.Range("N132").Value = -35030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9801.5
$ws.Range("I32").Value = 988.6786
$ws.Range("K32").Value = 988.6786
$ws.Range("M32").Value = -701.6786
$ws.Range("H97").Value = 1207.9375
$ws.Range("I97").Value = 760.75
$ws.Range("K97").Value = 760.75
$ws.Range("M97").Value = -264.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3066.1765
$ws.Range("J20").Value = 3912.8
$ws.Range("L20").Value = 3912.8
$ws.Range("N20").Value = -4406.8
$ws.Range("H86").Value = 4168741
$ws.Range("I86").Value = 5556975
$ws.Range("K86").Value = 5556975
$ws.Range("M86").Value = -5555852
$ws.Range("H89").Value = 4168741
$ws.Range("I89").Value = 5556975
$ws.Range("K89").Value = 27784875
$ws.Range("M89").Value = -27779259
$ws.Range("H107").Value = 2205.7778
$ws.Range("I107").Value = 1862.4615
$ws.Range("K107").Value = 1862.4615
$ws.Range("M107").Value = 57.53850000000011
$ws.Range("H108").Value = 100000
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 7787.9033
$ws.Range("I132").Value = 6865.125
$ws.Range("K132").Value = 20595.375
$ws.Range("M132").Value = -18065.375
$ws.Range("H134").Value = 11084.5
$ws.Range("I134").Value = 2221
$ws.Range("K134").Value = 6663
$ws.Range("M134").Value = -4128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 29617.143
$ws.Range("I46").Value = 1081.9524
$ws.Range("J46").Value = 72419.92999999999
$ws.Range("K46").Value = 3245.857199999999
$ws.Range("L46").Value = 217259.79
$ws.Range("M46").Value = -3154.857199999999
$ws.Range("N46").Value = -217441.79
$ws.Range("H61").Value = 402.57144
$ws.Range("J61").Value = 500
$ws.Range("L61").Value = 1500
$ws.Range("N61").Value = -1930
$ws.Range("H62").Value = 8665.666999999999
$ws.Range("I62").Value = 7498
$ws.Range("J62").Value = 9249.5
$ws.Range("K62").Value = 22494
$ws.Range("L62").Value = 27748.5
$ws.Range("M62").Value = -21808
$ws.Range("N62").Value = -29120.5
$ws.Range("H65").Value = 8665.666999999999
$ws.Range("I65").Value = 7498
$ws.Range("J65").Value = 9249.5
$ws.Range("K65").Value = 67482
$ws.Range("L65").Value = 83245.5
$ws.Range("M65").Value = -64050
$ws.Range("N65").Value = -90109.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5899.6665
$ws.Range("I70").Value = 6079.6
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 6079.6
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -5809.6
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 5899.6665
$ws.Range("I73").Value = 6079.6
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 6079.6
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -5143.6
$ws.Range("N73").Value = -6872
$ws.Range("H113").Value = 5794.9473
$ws.Range("I113").Value = 7168.1
$ws.Range("K113").Value = 7168.1
$ws.Range("M113").Value = -4998.1
$ws.Range("H126").Value = 8816
$ws.Range("I126").Value = 7999
$ws.Range("J126").Value = 9633
$ws.Range("K126").Value = 23997
$ws.Range("L126").Value = 28899
$ws.Range("M126").Value = -21527
$ws.Range("N126").Value = -33839
$ws.Range("H132").Value = 5470.129
$ws.Range("I132").Value = 5675.385
$ws.Range("J132").Value = 4402.8
$ws.Range("K132").Value = 17026.155
$ws.Range("L132").Value = 13208.4
$ws.Range("M132").Value = -14496.155
$ws.Range("N132").Value = -18268.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11206.186
$ws.Range("I7").Value = 11272.782
$ws.Range("K7").Value = 11272.782
$ws.Range("M7").Value = -11160.782
$ws.Range("H46").Value = 4990.35
$ws.Range("I46").Value = 6206.4287
$ws.Range("J46").Value = 2152.8333
$ws.Range("K46").Value = 6206.4287
$ws.Range("L46").Value = 2152.8333
$ws.Range("M46").Value = -6018.4287
$ws.Range("N46").Value = -2528.8333
$ws.Range("H61").Value = 2210.8667
$ws.Range("I61").Value = 2161.08
$ws.Range("J61").Value = 2459.8
$ws.Range("K61").Value = 2161.08
$ws.Range("L61").Value = 2459.8
$ws.Range("M61").Value = -1959.08
$ws.Range("N61").Value = -2863.8
$ws.Range("H113").Value = 2210.8667
$ws.Range("I113").Value = 2161.08
$ws.Range("J113").Value = 2459.8
$ws.Range("K113").Value = 2161.08
$ws.Range("L113").Value = 2459.8
$ws.Range("M113").Value = 8.920000000000073
$ws.Range("N113").Value = -6799.8
$ws.Range("H122").Value = 9055.1
$ws.Range("J122").Value = 3349.5
$ws.Range("L122").Value = 10048.5
$ws.Range("N122").Value = -14948.5
$ws.Range("H126").Value = 11206.186
$ws.Range("I126").Value = 11272.782
$ws.Range("K126").Value = 33818.346
$ws.Range("M126").Value = -31348.346

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4266.3477
$ws.Range("I132").Value = 4385.6313
$ws.Range("K132").Value = 13156.8939
$ws.Range("M132").Value = -10626.8939
